## Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": re-apply number/wrap formatting to rows 2-5 (A:F) and
# correct the tranche totals in row 4 (A4/B4 200 -> 100).
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("A2").Value = 10000
$wsSummary.Range("E2").Value = 10000
$wsSummary.Range("A2:F2").VerticalAlignment = -4108
$wsSummary.Range("A2:F2").WrapText = $true
$wsSummary.Range("A2").NumberFormat = "#,##0"
$wsSummary.Range("E2").NumberFormat = "#,##0"
$wsSummary.Range("B2").NumberFormat = "General"
$wsSummary.Range("C2").NumberFormat = "General"
$wsSummary.Range("D2").NumberFormat = "General"
$wsSummary.Range("F2").NumberFormat = "#,##0.00"

$wsSummary.Range("A3:F3").VerticalAlignment = -4108
$wsSummary.Range("A3:F3").WrapText = $true
$wsSummary.Range("A3:F3").NumberFormat = "General"

$wsSummary.Range("A4").Value = 100
$wsSummary.Range("B4").Value = 100
$wsSummary.Range("A4:F4").VerticalAlignment = -4108
$wsSummary.Range("A4:F4").WrapText = $true
$wsSummary.Range("A4:F4").NumberFormat = "General"

$wsSummary.Range("A5:F5").VerticalAlignment = -4108
$wsSummary.Range("A5:F5").WrapText = $true
$wsSummary.Range("A5:F5").NumberFormat = "General"

$wsSummary.Range("A7:XFD14").Select()

# ---------------------------------------------------------------------------
# Sheet "Repayment Schedule": the disbursement-fee tranche correction means
# row 4's fee columns (I/K/L) drop back to 0.
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Range("I4").Value = 0
$wsRepay.Range("K4").Value = 0
$wsRepay.Range("L4").Value = 0

$wsRepay.Range("M20:M21").Select()

# ---------------------------------------------------------------------------
# Sheet "Transactions": two new columns (K, L) were inserted, and the
# transaction rows were regenerated for the corrected Add-Tranche scenario.
# ---------------------------------------------------------------------------
$wsTx = $wb.Worksheets.Item("Transactions")

$wsTx.Columns("K:L").Insert()

$wsTx.Range("A2").Value = 13
$wsTx.Range("B2").Value = "Head Office"
$wsTx.Range("C2").Value = 42064
$wsTx.Range("D2").Value = "Disbursement"
$wsTx.Range("E2").Value = 5000
$wsTx.Range("F2").Value = 0
$wsTx.Range("G2").Value = 0
$wsTx.Range("H2").Value = 0
$wsTx.Range("I2").Value = 0
$wsTx.Range("J2").Value = 10000
$wsTx.Range("K2").Value = ""
$wsTx.Range("L2").Value = ""

$wsTx.Range("A3").Value = 11
$wsTx.Range("B3").Value = "Head Office"
$wsTx.Range("C3").Value = 42005
$wsTx.Range("D3").Value = "Repayment (at time of disbursement)"
$wsTx.Range("E3").Value = 100
$wsTx.Range("F3").Value = 0
$wsTx.Range("G3").Value = 0
$wsTx.Range("H3").Value = 100
$wsTx.Range("I3").Value = 0
$wsTx.Range("J3").Value = 5000
$wsTx.Range("K3").Value = ""
$wsTx.Range("L3").Value = ""

$wsTx.Range("A4").Value = 10
$wsTx.Range("B4").Value = "Head Office"
$wsTx.Range("C4").Value = 42005
$wsTx.Range("D4").Value = "Disbursement"
$wsTx.Range("E4").Value = 5000
$wsTx.Range("F4").Value = 0
$wsTx.Range("G4").Value = 0
$wsTx.Range("H4").Value = 0
$wsTx.Range("I4").Value = 0
$wsTx.Range("J4").Value = 5000

$wsTx.Range("A5:J5").ClearContents()

$wsTx.Range("A2:A4,B2:B4,D2:D4,F2:I4").NumberFormat = "General"
$wsTx.Range("A2:A4,B2:B4,D2:D4,F2:I4").HorizontalAlignment = -4131
$wsTx.Range("A2:A4,B2:B4,D2:D4,F2:I4").VerticalAlignment = -4160

$wsTx.Range("C2:C4").NumberFormat = "mm/dd/yyyy"
$wsTx.Range("C2:C4").HorizontalAlignment = -4131
$wsTx.Range("C2:C4").VerticalAlignment = -4160

$wsTx.Range("E2:E4,J2:J4").NumberFormat = "#,##0"
$wsTx.Range("E2:E4,J2:J4").HorizontalAlignment = -4131
$wsTx.Range("E2:E4,J2:J4").VerticalAlignment = -4160

$wsTx.Range("K2:L3").Font.Italic = $true
$wsTx.Range("K2:L3").VerticalAlignment = -4108
$wsTx.Range("K2:L3").WrapText = $true

$wsTx.Columns("A:A").ColumnWidth = 7.7109375
$wsTx.Columns("B:B").ColumnWidth = 10.5703125
$wsTx.Columns("C:C").ColumnWidth = 15.140625
$wsTx.Columns("D:D").ColumnWidth = 32.28515625
$wsTx.Columns("E:E").ColumnWidth = 7.42578125
$wsTx.Columns("F:F").ColumnWidth = 8.140625
$wsTx.Columns("G:G").ColumnWidth = 7.140625
$wsTx.Columns("H:H").ColumnWidth = 5.140625
$wsTx.Columns("I:I").ColumnWidth = 8.7109375
$wsTx.Columns("J:J").ColumnWidth = 12.28515625

$wsTx.Range("B3").Select()
